# Generate Report for Handoff
# Adds two new tracked files (8666816a-...md and e9e193ba-...md) to all
# three report sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$DATE_FMT = "yyyy-mm-dd HH:mm:ss"
$LINK_COLOR = 15570276   # BGR for RGB FF6495ED, matches the workbook's custom HyperLink style
$REPO = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f2d54e5d416216170d3d301d0ea81033c4d7af4/e2e/"

function Style-Link($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $LINK_COLOR
}

function Style-Date($rng) {
    $rng.NumberFormat = $DATE_FMT
}

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | Path And Name | Extension | Publish URL |
#                   zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @{ Row = 6; Name = "8666816a-aab3-472d-bea9-9e0ef1f71c72.md"; Path = "e2e\8666816a-aab3-472d-bea9-9e0ef1f71c72.md" },
    @{ Row = 7; Name = "e9e193ba-a693-4a31-8955-56b5a06a862f.md"; Path = "e2e\e9e193ba-a693-4a31-8955-56b5a06a862f.md" }
)

foreach ($item in $overviewRows) {
    $r = $item.Row
    $wsOverview.Range("A$r").Value = $item.Name
    $wsOverview.Range("B$r").Value = $item.Path
    $wsOverview.Hyperlinks.Add($wsOverview.Range("B$r"), ($REPO + $item.Name), "", "", $item.Path) | Out-Null
    Style-Link $wsOverview.Range("B$r")
    $wsOverview.Range("C$r").Value = ".md"
    $wsOverview.Range("D$r").Value = ""
    $wsOverview.Range("E$r").Value = "Ready for handoff"
    $wsOverview.Range("F$r").Value = "Ready for handoff"
    $wsOverview.Range("G$r").Value = "2016-08-19 02:39:06"
    Style-Date $wsOverview.Range("G$r")
}

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G7")) | Out-Null

# ---------------------------------------------------------------------
# Sheets "zh-cn" / "de-de": Source File Name | File Extension | Status |
#   Source Path | Priority | Content Duplicate | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Reference Tokens | To be localized |
#   Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------
function Fill-LocaleSheet($ws, $xlfSuffix, $hoDate) {
    $rows = @(
        @{ Row = 6; Name = "8666816a-aab3-472d-bea9-9e0ef1f71c72.md"; Xlf = "8666816a-aab3-472d-bea9-9e0ef1f71c72.181a07d97efc273f7fc03612a0027416064bd172.$xlfSuffix.xlf" },
        @{ Row = 7; Name = "e9e193ba-a693-4a31-8955-56b5a06a862f.md"; Xlf = "e9e193ba-a693-4a31-8955-56b5a06a862f.95ddb336a8441eeccdf4107c08cc0562383eda05.$xlfSuffix.xlf" }
    )

    foreach ($item in $rows) {
        $r = $item.Row
        $ws.Range("A$r").Value = $item.Name
        $ws.Hyperlinks.Add($ws.Range("A$r"), ($REPO + $item.Name), "", "", $item.Name) | Out-Null
        Style-Link $ws.Range("A$r")

        $ws.Range("B$r").Value = ".md"
        $ws.Range("C$r").Value = "Ready for handoff"
        $ws.Range("D$r").Value = "e2e"
        $ws.Range("E$r").Value = "'ht"
        $ws.Range("F$r").Value = "'False"
        $ws.Range("G$r").Value = $item.Xlf
        $ws.Range("H$r").Value = $hoDate
        Style-Date $ws.Range("H$r")
        $ws.Range("I$r").Value = ""
        $ws.Range("J$r").Value = ""
        $ws.Range("K$r").Value = "0001-01-01 00:00:00"
        Style-Date $ws.Range("K$r")
        $ws.Range("L$r").Value = ""
        $ws.Range("M$r").Value = "'True"
        $ws.Range("N$r").Value = ""
        $ws.Range("O$r").Value = "'False"
        $ws.Range("P$r").Value = ""
    }
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Fill-LocaleSheet $wsZhCn "zh-cn" "2016-08-19 02:38:57"
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P7")) | Out-Null

$wsDeDe = $wb.Worksheets.Item("de-de")
Fill-LocaleSheet $wsDeDe "de-de" "2016-08-19 02:39:06"
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P7")) | Out-Null

Write-Host "Handoff report rows added."
